$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 9 and Day 10 rows: only the "Day" column (A) is known so far.
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10

# Day 11: both the day number and its measured runtime are known.
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 0.00086514

# Leave the selection where the author ended up after entering today's data.
$ws.Range("B11").Select() | Out-Null
